$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Sending cluster" column (A) for rows 2-5: Resolving-Mac -> Inflammatory-Mac ---
$ws.Range("A2:A5").Value2 = "Inflammatory-Mac"

# --- Row 2 (target cluster: ECs) ---
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 0.7878926666666667
$ws.Range("H2").Value2 = 2.363678
$ws.Range("M2").Value2 = 1.834290333333333
$ws.Range("N2").Value2 = 5.502871
$ws.Range("O2").Value2 = 0.1177372815936135
$ws.Range("P2").Value2 = 0.1177372815936135
$ws.Range("Q2").Value2 = 1.445223902170889
$ws.Range("R2").Value2 = 13.007015119538
$ws.Range("S2").Value2 = 0.1177372815936135
$ws.Range("T2").Value2 = 0.1177372815936135

# --- Row 3 (target cluster: FAPs) ---
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 0.7878926666666667
$ws.Range("H3").Value2 = 2.363678
$ws.Range("O3").Value2 = 0.7242355134604062
$ws.Range("P3").Value2 = 0.7242355134604062
$ws.Range("Q3").Value2 = 8.889983365394444
$ws.Range("R3").Value2 = 80.00985028855
$ws.Range("S3").Value2 = 0.7242355134604062
$ws.Range("T3").Value2 = 0.7242355134604062

# --- Row 4 (target cluster: was Inflammatory-Mac, now MuSCs) ---
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 0.7878926666666667
$ws.Range("H4").Value2 = 2.363678
$ws.Range("M4").Value2 = 2.383963333333333
$ws.Range("N4").Value2 = 7.15189
$ws.Range("O4").Value2 = 0.1530190489394624
$ws.Range("P4").Value2 = 0.1530190489394624
$ws.Range("Q4").Value2 = 1.878307227935556
$ws.Range("R4").Value2 = 16.90476505142
$ws.Range("S4").Value2 = 0.1530190489394624
$ws.Range("T4").Value2 = 0.1530190489394624

# --- Row 5 (target cluster: was MuSCs, now Resolving-Mac) ---
$ws.Range("D5").Value2 = "Resolving-Mac"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 0.7878926666666667
$ws.Range("H5").Value2 = 2.363678
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 0.6666666666666666
$ws.Range("M5").Value2 = 0.07802466666666667
$ws.Range("N5").Value2 = 0.234074
$ws.Range("O5").Value2 = 0.005008156006517959
$ws.Range("P5").Value2 = 0.005008156006517959
$ws.Range("Q5").Value2 = 0.06147506268577779
$ws.Range("R5").Value2 = 0.5532755641720001
$ws.Range("S5").Value2 = 0.005008156006517959
$ws.Range("T5").Value2 = 0.005008156006517959

# --- Row 6: removed entirely (self-loop Resolving-Mac -> Resolving-Mac dropped) ---
$ws.Rows.Item(6).Delete()
